# Refresh cryptocurrency price/volume snapshot (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-like numbers (locale-formatted with dots as
# thousands separators, e.g. "62.507.57"). Force the cell to Text format
# before writing so Excel doesn't silently reinterpret the string as a
# floating point number / scientific notation, then restore the default
# "Normal" style so no stray number-format style is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.507.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.435.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.432.02'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("E10").Value = '  -3.82%  '
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("E15").Value = '  -3.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.874.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.255.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.442.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.70%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '630.14'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0961'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.38%  '
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.994'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("E31").Value = '  -3.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("E34").Value = '  -6.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.45'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.375'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.05%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.94'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.73'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '145.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.15%  '
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0524'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.597'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.73'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.32%  '
$ws.Range("E51").Value = '  -2.02%  '
